$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.164.26"
$ws.Range("E2").Value = "  +1.73%  "

$ws.Range("D3").Value = "2.524.47"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +4.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.41"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.26"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +9.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0822"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.92%  "

$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.84%  "

$ws.Range("D15").Value = "2.917.78"
$ws.Range("E15").Value = "  +0.73%  "

$ws.Range("D16").Value = "2.521.47"
$ws.Range("E16").Value = "  +0.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.862"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").Value = "48.072.10"
$ws.Range("E18").Value = "  +1.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.22%  "

$ws.Range("D21").Value = "0.0₃0948"
$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.72"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.42"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.88"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +9.56%  "

$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.14"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.03%  "

$ws.Range("E30").Value = "  +6.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.18"
$ws.Range("D31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.71"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.17%  "

$ws.Range("E34").Value = "  -0.95%  "

$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.73"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.39"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.99"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").Value = "2.007.20"
$ws.Range("E45").Value = "  +0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.16"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.89"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.17%  "

$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.11"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.09%  "
